$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Fransisco"
$ws.Range("B5").Value = 35
$ws.Range("C5").Value = 88
